$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A41").Value = "Michela Menghini"
$ws.Range("B41").Value = "Nicolas Giordani  | FC SAVIGNANO"
$ws.Range("C41").Value = "Andrea Conzatti | FC SAVIGNANO"
$ws.Range("D41").Value = "Matteo Mazzola | GREP"
$ws.Range("E41").Value = "Alessio Farinati | Pinguini Trentini"
$ws.Range("F41").Value = "Emanuele  valduga | wanda tim"
